# Apply updated crypto price/volume data (GitHub Actions scrape refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.859.03"
$ws.Range("E2").Value = "  -5.66%  "

$ws.Range("D3").Value = "1.819.93"
$ws.Range("E3").Value = "  -4.40%  "

$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  -0.36%  "

$ws.Range("D5").Value = "'328.50"
$ws.Range("E5").Value = "  -2.79%  "

$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "  -0.26%  "

$ws.Range("D7").Value = "'0.4639"
$ws.Range("E7").Value = "  -2.66%  "

$ws.Range("D8").Value = "'0.3847"
$ws.Range("E8").Value = "  -3.74%  "

$ws.Range("D9").Value = "'45.81"
$ws.Range("E9").Value = "  -3.06%  "

$ws.Range("D10").Value = "'0.07844"
$ws.Range("E10").Value = "  -2.44%  "

$ws.Range("D11").Value = "'0.9588"
$ws.Range("E11").Value = "  -3.20%  "

$ws.Range("D12").Value = "'21.80"
$ws.Range("E12").Value = "  -6.16%  "

$ws.Range("D13").Value = "1.828.53"
$ws.Range("E13").Value = "  -3.88%  "

$ws.Range("D14").Value = "'5.642"
$ws.Range("E14").Value = "  -4.52%  "

$ws.Range("D15").Value = "'6.849"
$ws.Range("E15").Value = "  -3.68%  "

$ws.Range("D16").Value = "'0.06858"
$ws.Range("E16").Value = "  +0.48%  "

$ws.Range("D17").Value = "'1.002"
$ws.Range("E17").Value = "  -0.37%  "

$ws.Range("D18").Value = "'86.61"
$ws.Range("E18").Value = "  -2.76%  "

$ws.Range("D19").Value = "'0.000009908"
$ws.Range("E19").Value = "  -2.91%  "

$ws.Range("D20").Value = "'16.61"
$ws.Range("E20").Value = "  -4.23%  "

$ws.Range("E21").Value = "  -0.41%  "

$ws.Range("D22").Value = "27.882.06"
$ws.Range("E22").Value = "  -5.63%  "

$ws.Range("D23").Value = "'5.309"
$ws.Range("E23").Value = "  -3.57%  "

$ws.Range("D24").Value = "'10.95"
$ws.Range("E24").Value = "  -5.63%  "

$ws.Range("D25").Value = "'2.088"
$ws.Range("E25").Value = "  -3.15%  "

$ws.Range("D26").Value = "2.041.26"
$ws.Range("E26").Value = "  -3.54%  "

$ws.Range("D27").Value = "'152.34"
$ws.Range("E27").Value = "  -2.79%  "

$ws.Range("D28").Value = "'19.19"
$ws.Range("E28").Value = "  -1.72%  "

$ws.Range("D29").Value = "'5.739"
$ws.Range("E29").Value = "  -11.81%  "

$ws.Range("D30").Value = "'1.968"
$ws.Range("E30").Value = "  -4.32%  "

$ws.Range("D31").Value = "'116.54"
$ws.Range("E31").Value = "  -2.12%  "

$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "'0.9339"
$ws.Range("E32").Value = "  -6.33%  "

$ws.Range("B33").Value = "Stellar"
$ws.Range("C33").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D33").Value = "'0.09238"
$ws.Range("E33").Value = "  -3.11%  "

$ws.Range("D34").Value = "'5.283"
$ws.Range("E34").Value = "  -3.33%  "

$ws.Range("D35").Value = "'3.351"
$ws.Range("E35").Value = "  -5.23%  "

$ws.Range("D36").Value = "'1.313"
$ws.Range("E36").Value = "  -5.30%  "

$ws.Range("D37").Value = "'0.05931"
$ws.Range("E37").Value = "  -8.43%  "

$ws.Range("D38").Value = "'0.02143"
$ws.Range("E38").Value = "  -4.35%  "

$ws.Range("D39").Value = "'1.143"
$ws.Range("E39").Value = "  -3.94%  "

$ws.Range("D40").Value = "'1.003"
$ws.Range("E40").Value = "  -0.27%  "

$ws.Range("D41").Value = "'7.566"
$ws.Range("E41").Value = "  -2.33%  "

$ws.Range("D42").Value = "'0.5568"
$ws.Range("E42").Value = "  -4.26%  "

$ws.Range("D43").Value = "'9.900"
$ws.Range("E43").Value = "  -5.96%  "

$ws.Range("D44").Value = "'0.1766"
$ws.Range("E44").Value = "  -2.86%  "

$ws.Range("D45").Value = "'1.211"
$ws.Range("E45").Value = "  -4.36%  "

$ws.Range("D46").Value = "'2.232"
$ws.Range("E46").Value = "  -8.76%  "

$ws.Range("D47").Value = "'11.54"
$ws.Range("E47").Value = "  -4.95%  "

$ws.Range("D48").Value = "'0.5235"
$ws.Range("E48").Value = "  -4.45%  "

$ws.Range("D49").Value = "'0.06981"
$ws.Range("E49").Value = "  -5.84%  "

$ws.Range("D50").Value = "'1.821"
$ws.Range("E50").Value = "  -6.66%  "

$ws.Range("D51").Value = "'112.34"
$ws.Range("E51").Value = "  -3.02%  "
